# Fixed column headers and input types
#
# The "Input Type" column used inconsistent, human-readable labels
# ("Dropdown", "Number input", "Text area") and the two trailing header
# cells had a typo ("optuons") / an unclear name ("Unit Options"). This
# normalizes the Input Type values to short machine-friendly tokens
# ("dropdown", "number", "text") and fixes the header row, plus a stray
# "ni" typo that had crept into the Unit Options column for row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) --------------------------------------------------
$ws.Cells.Item(1, 4).Value = "options"
$ws.Cells.Item(1, 5).Value = "unit_required"

# --- Fix the stray "ni" typo in row 3's unit/options column --------------
$ws.Cells.Item(3, 5).Value = "no"

# --- Normalize every "Input Type" value (column C, rows 2-38) ------------
$map = @{
    "Dropdown"      = "dropdown"
    "Number input"  = "number"
    "Text area"     = "text"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 38) { $lastRow = 38 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Text
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}

# --- Restore the view state (active cell / scroll position) --------------
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 1
$ws.Range("C41").Select()
